$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Ají" / "Cacho cabra verde" in the
# "Macroferia Regional de Talca" data set. It belongs right above the
# existing row 120, so insert a new row there and shift everything else
# (old rows 120..234) down by one, producing the new row 235 at the bottom.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new record. Most fields are
# identical to the neighbouring rows (same market/region/category/variety),
# only the date, volume and prices differ.
$ws.Cells.Item(120, 1).Value = 5
$ws.Cells.Item(120, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(120, 3).Value = "Maule"
$ws.Cells.Item(120, 4).Value = 44658
$ws.Cells.Item(120, 5).Value = 7
$ws.Cells.Item(120, 6).Value = 100112021
$ws.Cells.Item(120, 7).Value = "Ají"
$ws.Cells.Item(120, 8).Value = "Cacho cabra verde"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 100
$ws.Cells.Item(120, 11).Value = 14000
$ws.Cells.Item(120, 12).Value = 14000
$ws.Cells.Item(120, 13).Value = 14000
$ws.Cells.Item(120, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(120, 15).Value = "Región del Maule"
$ws.Cells.Item(120, 16).Value = 560
$ws.Cells.Item(120, 17).Value = 25
$ws.Cells.Item(120, 18).Value = "Hortaliza"
